$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate rows [19, 20, 21, 22, 23] (content shifts forward, last wraps to first) ---
${v_19_A} = 131024468
${v_19_B} = 57884
${v_19_D} = 'NT'
${v_19_E} = 100109
${v_19_F} = 'Tretåig hackspett'
${v_19_G} = 'Picoides tridactylus'
${v_19_H} = '(Linnaeus, 1758)'
${v_19_I} = ''
${v_19_K} = $null
${v_19_L} = $null
${v_19_M} = $null
${v_19_N} = $null
${v_19_P} = 'Sjöarsjön N, Jmt'
${v_19_Q} = 464332
${v_19_R} = 7042344
${v_19_S} = 10
${v_19_T} = 'Jämtland'
${v_19_U} = 'Krokom'
${v_19_V} = 'Jämtland'
${v_19_W} = 'Offerdal'
${v_19_Y} = '2026-02-03'
${v_19_AA} = '2026-02-03'
${v_19_AC} = 'Ringhack färska och äldre'
${v_19_AD} = $false
${v_19_AE} = $false
${v_19_AG} = $false
${v_19_AH} = $null
${v_19_AJ} = $null
${v_19_AK} = $null
${v_19_AM} = $null
${v_19_AO} = $null
${v_19_AT} = ''
${v_19_AW} = 'Benny Öwre'
${v_19_AX} = 'Benny Öwre'
${v_19_AY} = ''

${v_20_A} = 131024458
${v_20_B} = 57884
${v_20_D} = 'NT'
${v_20_E} = 100109
${v_20_F} = 'Tretåig hackspett'
${v_20_G} = 'Picoides tridactylus'
${v_20_H} = '(Linnaeus, 1758)'
${v_20_I} = ''
${v_20_K} = $null
${v_20_L} = $null
${v_20_M} = $null
${v_20_N} = $null
${v_20_P} = 'Sjöarsjön N, Jmt'
${v_20_Q} = 464405
${v_20_R} = 7042391
${v_20_S} = 10
${v_20_T} = 'Jämtland'
${v_20_U} = 'Krokom'
${v_20_V} = 'Jämtland'
${v_20_W} = 'Offerdal'
${v_20_Y} = '2026-02-03'
${v_20_AA} = '2026-02-03'
${v_20_AC} = 'Ringhack'
${v_20_AD} = $false
${v_20_AE} = $false
${v_20_AG} = $false
${v_20_AH} = $null
${v_20_AJ} = $null
${v_20_AK} = $null
${v_20_AM} = $null
${v_20_AO} = $null
${v_20_AT} = ''
${v_20_AW} = 'Benny Öwre'
${v_20_AX} = 'Benny Öwre'
${v_20_AY} = ''

${v_21_A} = 131024464
${v_21_B} = 57884
${v_21_D} = 'NT'
${v_21_E} = 100109
${v_21_F} = 'Tretåig hackspett'
${v_21_G} = 'Picoides tridactylus'
${v_21_H} = '(Linnaeus, 1758)'
${v_21_I} = ''
${v_21_K} = $null
${v_21_L} = $null
${v_21_M} = $null
${v_21_N} = $null
${v_21_P} = 'Sjöarsjön N, Jmt'
${v_21_Q} = 464337
${v_21_R} = 7042380
${v_21_S} = 10
${v_21_T} = 'Jämtland'
${v_21_U} = 'Krokom'
${v_21_V} = 'Jämtland'
${v_21_W} = 'Offerdal'
${v_21_Y} = '2026-02-03'
${v_21_AA} = '2026-02-03'
${v_21_AC} = 'Ringhack färska och äldre'
${v_21_AD} = $false
${v_21_AE} = $false
${v_21_AG} = $false
${v_21_AH} = $null
${v_21_AJ} = $null
${v_21_AK} = $null
${v_21_AM} = $null
${v_21_AO} = $null
${v_21_AT} = ''
${v_21_AW} = 'Benny Öwre'
${v_21_AX} = 'Benny Öwre'
${v_21_AY} = ''

${v_22_A} = 131025468
${v_22_B} = 57884
${v_22_D} = 'NT'
${v_22_E} = 100109
${v_22_F} = 'Tretåig hackspett'
${v_22_G} = 'Picoides tridactylus'
${v_22_H} = '(Linnaeus, 1758)'
${v_22_I} = ''
${v_22_K} = ''
${v_22_L} = ''
${v_22_M} = 'äldre spår'
${v_22_N} = ''
${v_22_P} = 'Sjöarsjön Nordöst, Jmt'
${v_22_Q} = 464445
${v_22_R} = 7042578
${v_22_S} = 10
${v_22_T} = 'Jämtland'
${v_22_U} = 'Krokom'
${v_22_V} = 'Jämtland'
${v_22_W} = 'Offerdal'
${v_22_Y} = '2026-02-03'
${v_22_AA} = '2026-02-03'
${v_22_AC} = 'Ringhack, äldre, på en gran.'
${v_22_AD} = $false
${v_22_AE} = $false
${v_22_AG} = $false
${v_22_AH} = 'Granskog'
${v_22_AJ} = 'gran'
${v_22_AK} = 'Picea abies'
${v_22_AM} = 'Trädstam på levande träd'
${v_22_AO} = 'Stem on living tree # Picea abies'
${v_22_AT} = ''
${v_22_AW} = 'Kristian Zackrisson'
${v_22_AX} = 'Kristian Zackrisson'
${v_22_AY} = ''

${v_23_A} = 131025525
${v_23_B} = 79243
${v_23_D} = 'NT'
${v_23_E} = 6425
${v_23_F} = 'Garnlav'
${v_23_G} = 'Alectoria sarmentosa'
${v_23_H} = '(Ach.) Ach.'
${v_23_I} = ''
${v_23_K} = $null
${v_23_L} = $null
${v_23_M} = $null
${v_23_N} = $null
${v_23_P} = 'Sjöarsjön Nordöst, Jmt'
${v_23_Q} = 464067
${v_23_R} = 7042585
${v_23_S} = 10
${v_23_T} = 'Jämtland'
${v_23_U} = 'Krokom'
${v_23_V} = 'Jämtland'
${v_23_W} = 'Offerdal'
${v_23_Y} = '2026-02-03'
${v_23_AA} = '2026-02-03'
${v_23_AC} = 'På gran.'
${v_23_AD} = $false
${v_23_AE} = $false
${v_23_AG} = $false
${v_23_AH} = $null
${v_23_AJ} = $null
${v_23_AK} = $null
${v_23_AM} = $null
${v_23_AO} = $null
${v_23_AT} = ''
${v_23_AW} = 'Kristian Zackrisson'
${v_23_AX} = 'Kristian Zackrisson'
${v_23_AY} = ''

$ws.Cells.Item(20, 1).Value2 = ${v_19_A}
$ws.Cells.Item(20, 2).Value2 = ${v_19_B}
$ws.Cells.Item(20, 4).Value2 = ${v_19_D}
$ws.Cells.Item(20, 5).Value2 = ${v_19_E}
$ws.Cells.Item(20, 6).Value2 = ${v_19_F}
$ws.Cells.Item(20, 7).Value2 = ${v_19_G}
$ws.Cells.Item(20, 8).Value2 = ${v_19_H}
$ws.Cells.Item(20, 9).Value2 = ${v_19_I}
$ws.Cells.Item(20, 11).Value2 = ${v_19_K}
$ws.Cells.Item(20, 12).Value2 = ${v_19_L}
$ws.Cells.Item(20, 13).Value2 = ${v_19_M}
$ws.Cells.Item(20, 14).Value2 = ${v_19_N}
$ws.Cells.Item(20, 16).Value2 = ${v_19_P}
$ws.Cells.Item(20, 17).Value2 = ${v_19_Q}
$ws.Cells.Item(20, 18).Value2 = ${v_19_R}
$ws.Cells.Item(20, 19).Value2 = ${v_19_S}
$ws.Cells.Item(20, 20).Value2 = ${v_19_T}
$ws.Cells.Item(20, 21).Value2 = ${v_19_U}
$ws.Cells.Item(20, 22).Value2 = ${v_19_V}
$ws.Cells.Item(20, 23).Value2 = ${v_19_W}
$ws.Cells.Item(20, 25).NumberFormat = "@"
$ws.Cells.Item(20, 25).Value2 = ${v_19_Y}
$ws.Cells.Item(20, 27).NumberFormat = "@"
$ws.Cells.Item(20, 27).Value2 = ${v_19_AA}
$ws.Cells.Item(20, 29).Value2 = ${v_19_AC}
$ws.Cells.Item(20, 30).Value2 = ${v_19_AD}
$ws.Cells.Item(20, 31).Value2 = ${v_19_AE}
$ws.Cells.Item(20, 33).Value2 = ${v_19_AG}
$ws.Cells.Item(20, 34).Value2 = ${v_19_AH}
$ws.Cells.Item(20, 36).Value2 = ${v_19_AJ}
$ws.Cells.Item(20, 37).Value2 = ${v_19_AK}
$ws.Cells.Item(20, 39).Value2 = ${v_19_AM}
$ws.Cells.Item(20, 41).Value2 = ${v_19_AO}
$ws.Cells.Item(20, 46).Value2 = ${v_19_AT}
$ws.Cells.Item(20, 49).Value2 = ${v_19_AW}
$ws.Cells.Item(20, 50).Value2 = ${v_19_AX}
$ws.Cells.Item(20, 51).Value2 = ${v_19_AY}

$ws.Cells.Item(21, 1).Value2 = ${v_20_A}
$ws.Cells.Item(21, 2).Value2 = ${v_20_B}
$ws.Cells.Item(21, 4).Value2 = ${v_20_D}
$ws.Cells.Item(21, 5).Value2 = ${v_20_E}
$ws.Cells.Item(21, 6).Value2 = ${v_20_F}
$ws.Cells.Item(21, 7).Value2 = ${v_20_G}
$ws.Cells.Item(21, 8).Value2 = ${v_20_H}
$ws.Cells.Item(21, 9).Value2 = ${v_20_I}
$ws.Cells.Item(21, 11).Value2 = ${v_20_K}
$ws.Cells.Item(21, 12).Value2 = ${v_20_L}
$ws.Cells.Item(21, 13).Value2 = ${v_20_M}
$ws.Cells.Item(21, 14).Value2 = ${v_20_N}
$ws.Cells.Item(21, 16).Value2 = ${v_20_P}
$ws.Cells.Item(21, 17).Value2 = ${v_20_Q}
$ws.Cells.Item(21, 18).Value2 = ${v_20_R}
$ws.Cells.Item(21, 19).Value2 = ${v_20_S}
$ws.Cells.Item(21, 20).Value2 = ${v_20_T}
$ws.Cells.Item(21, 21).Value2 = ${v_20_U}
$ws.Cells.Item(21, 22).Value2 = ${v_20_V}
$ws.Cells.Item(21, 23).Value2 = ${v_20_W}
$ws.Cells.Item(21, 25).NumberFormat = "@"
$ws.Cells.Item(21, 25).Value2 = ${v_20_Y}
$ws.Cells.Item(21, 27).NumberFormat = "@"
$ws.Cells.Item(21, 27).Value2 = ${v_20_AA}
$ws.Cells.Item(21, 29).Value2 = ${v_20_AC}
$ws.Cells.Item(21, 30).Value2 = ${v_20_AD}
$ws.Cells.Item(21, 31).Value2 = ${v_20_AE}
$ws.Cells.Item(21, 33).Value2 = ${v_20_AG}
$ws.Cells.Item(21, 34).Value2 = ${v_20_AH}
$ws.Cells.Item(21, 36).Value2 = ${v_20_AJ}
$ws.Cells.Item(21, 37).Value2 = ${v_20_AK}
$ws.Cells.Item(21, 39).Value2 = ${v_20_AM}
$ws.Cells.Item(21, 41).Value2 = ${v_20_AO}
$ws.Cells.Item(21, 46).Value2 = ${v_20_AT}
$ws.Cells.Item(21, 49).Value2 = ${v_20_AW}
$ws.Cells.Item(21, 50).Value2 = ${v_20_AX}
$ws.Cells.Item(21, 51).Value2 = ${v_20_AY}

$ws.Cells.Item(22, 1).Value2 = ${v_21_A}
$ws.Cells.Item(22, 2).Value2 = ${v_21_B}
$ws.Cells.Item(22, 4).Value2 = ${v_21_D}
$ws.Cells.Item(22, 5).Value2 = ${v_21_E}
$ws.Cells.Item(22, 6).Value2 = ${v_21_F}
$ws.Cells.Item(22, 7).Value2 = ${v_21_G}
$ws.Cells.Item(22, 8).Value2 = ${v_21_H}
$ws.Cells.Item(22, 9).Value2 = ${v_21_I}
$ws.Cells.Item(22, 11).Value2 = ${v_21_K}
$ws.Cells.Item(22, 12).Value2 = ${v_21_L}
$ws.Cells.Item(22, 13).Value2 = ${v_21_M}
$ws.Cells.Item(22, 14).Value2 = ${v_21_N}
$ws.Cells.Item(22, 16).Value2 = ${v_21_P}
$ws.Cells.Item(22, 17).Value2 = ${v_21_Q}
$ws.Cells.Item(22, 18).Value2 = ${v_21_R}
$ws.Cells.Item(22, 19).Value2 = ${v_21_S}
$ws.Cells.Item(22, 20).Value2 = ${v_21_T}
$ws.Cells.Item(22, 21).Value2 = ${v_21_U}
$ws.Cells.Item(22, 22).Value2 = ${v_21_V}
$ws.Cells.Item(22, 23).Value2 = ${v_21_W}
$ws.Cells.Item(22, 25).NumberFormat = "@"
$ws.Cells.Item(22, 25).Value2 = ${v_21_Y}
$ws.Cells.Item(22, 27).NumberFormat = "@"
$ws.Cells.Item(22, 27).Value2 = ${v_21_AA}
$ws.Cells.Item(22, 29).Value2 = ${v_21_AC}
$ws.Cells.Item(22, 30).Value2 = ${v_21_AD}
$ws.Cells.Item(22, 31).Value2 = ${v_21_AE}
$ws.Cells.Item(22, 33).Value2 = ${v_21_AG}
$ws.Cells.Item(22, 34).Value2 = ${v_21_AH}
$ws.Cells.Item(22, 36).Value2 = ${v_21_AJ}
$ws.Cells.Item(22, 37).Value2 = ${v_21_AK}
$ws.Cells.Item(22, 39).Value2 = ${v_21_AM}
$ws.Cells.Item(22, 41).Value2 = ${v_21_AO}
$ws.Cells.Item(22, 46).Value2 = ${v_21_AT}
$ws.Cells.Item(22, 49).Value2 = ${v_21_AW}
$ws.Cells.Item(22, 50).Value2 = ${v_21_AX}
$ws.Cells.Item(22, 51).Value2 = ${v_21_AY}

$ws.Cells.Item(23, 1).Value2 = ${v_22_A}
$ws.Cells.Item(23, 2).Value2 = ${v_22_B}
$ws.Cells.Item(23, 4).Value2 = ${v_22_D}
$ws.Cells.Item(23, 5).Value2 = ${v_22_E}
$ws.Cells.Item(23, 6).Value2 = ${v_22_F}
$ws.Cells.Item(23, 7).Value2 = ${v_22_G}
$ws.Cells.Item(23, 8).Value2 = ${v_22_H}
$ws.Cells.Item(23, 9).Value2 = ${v_22_I}
$ws.Cells.Item(23, 11).Value2 = ${v_22_K}
$ws.Cells.Item(23, 12).Value2 = ${v_22_L}
$ws.Cells.Item(23, 13).Value2 = ${v_22_M}
$ws.Cells.Item(23, 14).Value2 = ${v_22_N}
$ws.Cells.Item(23, 16).Value2 = ${v_22_P}
$ws.Cells.Item(23, 17).Value2 = ${v_22_Q}
$ws.Cells.Item(23, 18).Value2 = ${v_22_R}
$ws.Cells.Item(23, 19).Value2 = ${v_22_S}
$ws.Cells.Item(23, 20).Value2 = ${v_22_T}
$ws.Cells.Item(23, 21).Value2 = ${v_22_U}
$ws.Cells.Item(23, 22).Value2 = ${v_22_V}
$ws.Cells.Item(23, 23).Value2 = ${v_22_W}
$ws.Cells.Item(23, 25).NumberFormat = "@"
$ws.Cells.Item(23, 25).Value2 = ${v_22_Y}
$ws.Cells.Item(23, 27).NumberFormat = "@"
$ws.Cells.Item(23, 27).Value2 = ${v_22_AA}
$ws.Cells.Item(23, 29).Value2 = ${v_22_AC}
$ws.Cells.Item(23, 30).Value2 = ${v_22_AD}
$ws.Cells.Item(23, 31).Value2 = ${v_22_AE}
$ws.Cells.Item(23, 33).Value2 = ${v_22_AG}
$ws.Cells.Item(23, 34).Value2 = ${v_22_AH}
$ws.Cells.Item(23, 36).Value2 = ${v_22_AJ}
$ws.Cells.Item(23, 37).Value2 = ${v_22_AK}
$ws.Cells.Item(23, 39).Value2 = ${v_22_AM}
$ws.Cells.Item(23, 41).Value2 = ${v_22_AO}
$ws.Cells.Item(23, 46).Value2 = ${v_22_AT}
$ws.Cells.Item(23, 49).Value2 = ${v_22_AW}
$ws.Cells.Item(23, 50).Value2 = ${v_22_AX}
$ws.Cells.Item(23, 51).Value2 = ${v_22_AY}

$ws.Cells.Item(19, 1).Value2 = ${v_23_A}
$ws.Cells.Item(19, 2).Value2 = ${v_23_B}
$ws.Cells.Item(19, 4).Value2 = ${v_23_D}
$ws.Cells.Item(19, 5).Value2 = ${v_23_E}
$ws.Cells.Item(19, 6).Value2 = ${v_23_F}
$ws.Cells.Item(19, 7).Value2 = ${v_23_G}
$ws.Cells.Item(19, 8).Value2 = ${v_23_H}
$ws.Cells.Item(19, 9).Value2 = ${v_23_I}
$ws.Cells.Item(19, 11).Value2 = ${v_23_K}
$ws.Cells.Item(19, 12).Value2 = ${v_23_L}
$ws.Cells.Item(19, 13).Value2 = ${v_23_M}
$ws.Cells.Item(19, 14).Value2 = ${v_23_N}
$ws.Cells.Item(19, 16).Value2 = ${v_23_P}
$ws.Cells.Item(19, 17).Value2 = ${v_23_Q}
$ws.Cells.Item(19, 18).Value2 = ${v_23_R}
$ws.Cells.Item(19, 19).Value2 = ${v_23_S}
$ws.Cells.Item(19, 20).Value2 = ${v_23_T}
$ws.Cells.Item(19, 21).Value2 = ${v_23_U}
$ws.Cells.Item(19, 22).Value2 = ${v_23_V}
$ws.Cells.Item(19, 23).Value2 = ${v_23_W}
$ws.Cells.Item(19, 25).NumberFormat = "@"
$ws.Cells.Item(19, 25).Value2 = ${v_23_Y}
$ws.Cells.Item(19, 27).NumberFormat = "@"
$ws.Cells.Item(19, 27).Value2 = ${v_23_AA}
$ws.Cells.Item(19, 29).Value2 = ${v_23_AC}
$ws.Cells.Item(19, 30).Value2 = ${v_23_AD}
$ws.Cells.Item(19, 31).Value2 = ${v_23_AE}
$ws.Cells.Item(19, 33).Value2 = ${v_23_AG}
$ws.Cells.Item(19, 34).Value2 = ${v_23_AH}
$ws.Cells.Item(19, 36).Value2 = ${v_23_AJ}
$ws.Cells.Item(19, 37).Value2 = ${v_23_AK}
$ws.Cells.Item(19, 39).Value2 = ${v_23_AM}
$ws.Cells.Item(19, 41).Value2 = ${v_23_AO}
$ws.Cells.Item(19, 46).Value2 = ${v_23_AT}
$ws.Cells.Item(19, 49).Value2 = ${v_23_AW}
$ws.Cells.Item(19, 50).Value2 = ${v_23_AX}
$ws.Cells.Item(19, 51).Value2 = ${v_23_AY}

# --- Rotate rows [33, 35, 34] (content shifts forward, last wraps to first) ---
${v_33_A} = 131024446
${v_33_B} = 57884
${v_33_D} = 'NT'
${v_33_E} = 100109
${v_33_F} = 'Tretåig hackspett'
${v_33_G} = 'Picoides tridactylus'
${v_33_H} = '(Linnaeus, 1758)'
${v_33_I} = ''
${v_33_K} = ''
${v_33_L} = ''
${v_33_M} = 'gammalt bo'
${v_33_N} = ''
${v_33_P} = 'Sjöarsjön N, Jmt'
${v_33_Q} = 463828
${v_33_R} = 7042514
${v_33_S} = 10
${v_33_T} = 'Jämtland'
${v_33_U} = 'Krokom'
${v_33_V} = 'Jämtland'
${v_33_W} = 'Offerdal'
${v_33_Y} = '2026-02-03'
${v_33_AA} = '2026-02-03'
${v_33_AC} = 'Gammalt bohål i grantickerötad gran'
${v_33_AD} = $false
${v_33_AE} = $true
${v_33_AG} = $false
${v_33_AH} = $null
${v_33_AJ} = $null
${v_33_AK} = $null
${v_33_AM} = $null
${v_33_AO} = $null
${v_33_AT} = ''
${v_33_AW} = 'Benny Öwre'
${v_33_AX} = 'Benny Öwre'
${v_33_AY} = ''

${v_35_A} = 131025464
${v_35_B} = 57884
${v_35_D} = 'NT'
${v_35_E} = 100109
${v_35_F} = 'Tretåig hackspett'
${v_35_G} = 'Picoides tridactylus'
${v_35_H} = '(Linnaeus, 1758)'
${v_35_I} = ''
${v_35_K} = ''
${v_35_L} = ''
${v_35_M} = 'färska spår'
${v_35_N} = ''
${v_35_P} = 'Sjöarsjön Nordöst, Jmt'
${v_35_Q} = 464210
${v_35_R} = 7042574
${v_35_S} = 10
${v_35_T} = 'Jämtland'
${v_35_U} = 'Krokom'
${v_35_V} = 'Jämtland'
${v_35_W} = 'Offerdal'
${v_35_Y} = '2026-02-03'
${v_35_AA} = '2026-02-03'
${v_35_AC} = 'Ringhack, färska och äldre, på stambasen av en gran.'
${v_35_AD} = $false
${v_35_AE} = $false
${v_35_AG} = $false
${v_35_AH} = 'Barrskog'
${v_35_AJ} = 'gran'
${v_35_AK} = 'Picea abies'
${v_35_AM} = 'Trädstam på levande träd'
${v_35_AO} = 'Stem on living tree # Picea abies'
${v_35_AT} = ''
${v_35_AW} = 'Kristian Zackrisson'
${v_35_AX} = 'Kristian Zackrisson'
${v_35_AY} = ''

${v_34_A} = 131024465
${v_34_B} = 57884
${v_34_D} = 'NT'
${v_34_E} = 100109
${v_34_F} = 'Tretåig hackspett'
${v_34_G} = 'Picoides tridactylus'
${v_34_H} = '(Linnaeus, 1758)'
${v_34_I} = ''
${v_34_K} = $null
${v_34_L} = $null
${v_34_M} = $null
${v_34_N} = $null
${v_34_P} = 'Sjöarsjön N, Jmt'
${v_34_Q} = 464328
${v_34_R} = 7042367
${v_34_S} = 10
${v_34_T} = 'Jämtland'
${v_34_U} = 'Krokom'
${v_34_V} = 'Jämtland'
${v_34_W} = 'Offerdal'
${v_34_Y} = '2026-02-03'
${v_34_AA} = '2026-02-03'
${v_34_AC} = 'Ringhack färska och äldre'
${v_34_AD} = $false
${v_34_AE} = $false
${v_34_AG} = $false
${v_34_AH} = $null
${v_34_AJ} = $null
${v_34_AK} = $null
${v_34_AM} = $null
${v_34_AO} = $null
${v_34_AT} = ''
${v_34_AW} = 'Benny Öwre'
${v_34_AX} = 'Benny Öwre'
${v_34_AY} = ''

$ws.Cells.Item(35, 1).Value2 = ${v_33_A}
$ws.Cells.Item(35, 2).Value2 = ${v_33_B}
$ws.Cells.Item(35, 4).Value2 = ${v_33_D}
$ws.Cells.Item(35, 5).Value2 = ${v_33_E}
$ws.Cells.Item(35, 6).Value2 = ${v_33_F}
$ws.Cells.Item(35, 7).Value2 = ${v_33_G}
$ws.Cells.Item(35, 8).Value2 = ${v_33_H}
$ws.Cells.Item(35, 9).Value2 = ${v_33_I}
$ws.Cells.Item(35, 11).Value2 = ${v_33_K}
$ws.Cells.Item(35, 12).Value2 = ${v_33_L}
$ws.Cells.Item(35, 13).Value2 = ${v_33_M}
$ws.Cells.Item(35, 14).Value2 = ${v_33_N}
$ws.Cells.Item(35, 16).Value2 = ${v_33_P}
$ws.Cells.Item(35, 17).Value2 = ${v_33_Q}
$ws.Cells.Item(35, 18).Value2 = ${v_33_R}
$ws.Cells.Item(35, 19).Value2 = ${v_33_S}
$ws.Cells.Item(35, 20).Value2 = ${v_33_T}
$ws.Cells.Item(35, 21).Value2 = ${v_33_U}
$ws.Cells.Item(35, 22).Value2 = ${v_33_V}
$ws.Cells.Item(35, 23).Value2 = ${v_33_W}
$ws.Cells.Item(35, 25).NumberFormat = "@"
$ws.Cells.Item(35, 25).Value2 = ${v_33_Y}
$ws.Cells.Item(35, 27).NumberFormat = "@"
$ws.Cells.Item(35, 27).Value2 = ${v_33_AA}
$ws.Cells.Item(35, 29).Value2 = ${v_33_AC}
$ws.Cells.Item(35, 30).Value2 = ${v_33_AD}
$ws.Cells.Item(35, 31).Value2 = ${v_33_AE}
$ws.Cells.Item(35, 33).Value2 = ${v_33_AG}
$ws.Cells.Item(35, 34).Value2 = ${v_33_AH}
$ws.Cells.Item(35, 36).Value2 = ${v_33_AJ}
$ws.Cells.Item(35, 37).Value2 = ${v_33_AK}
$ws.Cells.Item(35, 39).Value2 = ${v_33_AM}
$ws.Cells.Item(35, 41).Value2 = ${v_33_AO}
$ws.Cells.Item(35, 46).Value2 = ${v_33_AT}
$ws.Cells.Item(35, 49).Value2 = ${v_33_AW}
$ws.Cells.Item(35, 50).Value2 = ${v_33_AX}
$ws.Cells.Item(35, 51).Value2 = ${v_33_AY}

$ws.Cells.Item(34, 1).Value2 = ${v_35_A}
$ws.Cells.Item(34, 2).Value2 = ${v_35_B}
$ws.Cells.Item(34, 4).Value2 = ${v_35_D}
$ws.Cells.Item(34, 5).Value2 = ${v_35_E}
$ws.Cells.Item(34, 6).Value2 = ${v_35_F}
$ws.Cells.Item(34, 7).Value2 = ${v_35_G}
$ws.Cells.Item(34, 8).Value2 = ${v_35_H}
$ws.Cells.Item(34, 9).Value2 = ${v_35_I}
$ws.Cells.Item(34, 11).Value2 = ${v_35_K}
$ws.Cells.Item(34, 12).Value2 = ${v_35_L}
$ws.Cells.Item(34, 13).Value2 = ${v_35_M}
$ws.Cells.Item(34, 14).Value2 = ${v_35_N}
$ws.Cells.Item(34, 16).Value2 = ${v_35_P}
$ws.Cells.Item(34, 17).Value2 = ${v_35_Q}
$ws.Cells.Item(34, 18).Value2 = ${v_35_R}
$ws.Cells.Item(34, 19).Value2 = ${v_35_S}
$ws.Cells.Item(34, 20).Value2 = ${v_35_T}
$ws.Cells.Item(34, 21).Value2 = ${v_35_U}
$ws.Cells.Item(34, 22).Value2 = ${v_35_V}
$ws.Cells.Item(34, 23).Value2 = ${v_35_W}
$ws.Cells.Item(34, 25).NumberFormat = "@"
$ws.Cells.Item(34, 25).Value2 = ${v_35_Y}
$ws.Cells.Item(34, 27).NumberFormat = "@"
$ws.Cells.Item(34, 27).Value2 = ${v_35_AA}
$ws.Cells.Item(34, 29).Value2 = ${v_35_AC}
$ws.Cells.Item(34, 30).Value2 = ${v_35_AD}
$ws.Cells.Item(34, 31).Value2 = ${v_35_AE}
$ws.Cells.Item(34, 33).Value2 = ${v_35_AG}
$ws.Cells.Item(34, 34).Value2 = ${v_35_AH}
$ws.Cells.Item(34, 36).Value2 = ${v_35_AJ}
$ws.Cells.Item(34, 37).Value2 = ${v_35_AK}
$ws.Cells.Item(34, 39).Value2 = ${v_35_AM}
$ws.Cells.Item(34, 41).Value2 = ${v_35_AO}
$ws.Cells.Item(34, 46).Value2 = ${v_35_AT}
$ws.Cells.Item(34, 49).Value2 = ${v_35_AW}
$ws.Cells.Item(34, 50).Value2 = ${v_35_AX}
$ws.Cells.Item(34, 51).Value2 = ${v_35_AY}

$ws.Cells.Item(33, 1).Value2 = ${v_34_A}
$ws.Cells.Item(33, 2).Value2 = ${v_34_B}
$ws.Cells.Item(33, 4).Value2 = ${v_34_D}
$ws.Cells.Item(33, 5).Value2 = ${v_34_E}
$ws.Cells.Item(33, 6).Value2 = ${v_34_F}
$ws.Cells.Item(33, 7).Value2 = ${v_34_G}
$ws.Cells.Item(33, 8).Value2 = ${v_34_H}
$ws.Cells.Item(33, 9).Value2 = ${v_34_I}
$ws.Cells.Item(33, 11).Value2 = ${v_34_K}
$ws.Cells.Item(33, 12).Value2 = ${v_34_L}
$ws.Cells.Item(33, 13).Value2 = ${v_34_M}
$ws.Cells.Item(33, 14).Value2 = ${v_34_N}
$ws.Cells.Item(33, 16).Value2 = ${v_34_P}
$ws.Cells.Item(33, 17).Value2 = ${v_34_Q}
$ws.Cells.Item(33, 18).Value2 = ${v_34_R}
$ws.Cells.Item(33, 19).Value2 = ${v_34_S}
$ws.Cells.Item(33, 20).Value2 = ${v_34_T}
$ws.Cells.Item(33, 21).Value2 = ${v_34_U}
$ws.Cells.Item(33, 22).Value2 = ${v_34_V}
$ws.Cells.Item(33, 23).Value2 = ${v_34_W}
$ws.Cells.Item(33, 25).NumberFormat = "@"
$ws.Cells.Item(33, 25).Value2 = ${v_34_Y}
$ws.Cells.Item(33, 27).NumberFormat = "@"
$ws.Cells.Item(33, 27).Value2 = ${v_34_AA}
$ws.Cells.Item(33, 29).Value2 = ${v_34_AC}
$ws.Cells.Item(33, 30).Value2 = ${v_34_AD}
$ws.Cells.Item(33, 31).Value2 = ${v_34_AE}
$ws.Cells.Item(33, 33).Value2 = ${v_34_AG}
$ws.Cells.Item(33, 34).Value2 = ${v_34_AH}
$ws.Cells.Item(33, 36).Value2 = ${v_34_AJ}
$ws.Cells.Item(33, 37).Value2 = ${v_34_AK}
$ws.Cells.Item(33, 39).Value2 = ${v_34_AM}
$ws.Cells.Item(33, 41).Value2 = ${v_34_AO}
$ws.Cells.Item(33, 46).Value2 = ${v_34_AT}
$ws.Cells.Item(33, 49).Value2 = ${v_34_AW}
$ws.Cells.Item(33, 50).Value2 = ${v_34_AX}
$ws.Cells.Item(33, 51).Value2 = ${v_34_AY}

# --- Rotate rows [63, 64] (content shifts forward, last wraps to first) ---
${v_63_A} = 131024460
${v_63_B} = 57884
${v_63_D} = 'NT'
${v_63_E} = 100109
${v_63_F} = 'Tretåig hackspett'
${v_63_G} = 'Picoides tridactylus'
${v_63_H} = '(Linnaeus, 1758)'
${v_63_I} = ''
${v_63_P} = 'Sjöarsjön N, Jmt'
${v_63_Q} = 464385
${v_63_R} = 7042401
${v_63_S} = 10
${v_63_T} = 'Jämtland'
${v_63_U} = 'Krokom'
${v_63_V} = 'Jämtland'
${v_63_W} = 'Offerdal'
${v_63_Y} = '2026-02-03'
${v_63_AA} = '2026-02-03'
${v_63_AC} = 'Ringhack äldre'
${v_63_AD} = $false
${v_63_AE} = $false
${v_63_AG} = $false
${v_63_AT} = ''
${v_63_AW} = 'Benny Öwre'
${v_63_AX} = 'Benny Öwre'
${v_63_AY} = ''

${v_64_A} = 131024450
${v_64_B} = 57884
${v_64_D} = 'NT'
${v_64_E} = 100109
${v_64_F} = 'Tretåig hackspett'
${v_64_G} = 'Picoides tridactylus'
${v_64_H} = '(Linnaeus, 1758)'
${v_64_I} = ''
${v_64_P} = 'Sjöarsjön N, Jmt'
${v_64_Q} = 464273
${v_64_R} = 7042500
${v_64_S} = 10
${v_64_T} = 'Jämtland'
${v_64_U} = 'Krokom'
${v_64_V} = 'Jämtland'
${v_64_W} = 'Offerdal'
${v_64_Y} = '2026-02-03'
${v_64_AA} = '2026-02-03'
${v_64_AC} = 'Ringhack äldre'
${v_64_AD} = $false
${v_64_AE} = $false
${v_64_AG} = $false
${v_64_AT} = ''
${v_64_AW} = 'Benny Öwre'
${v_64_AX} = 'Benny Öwre'
${v_64_AY} = ''

$ws.Cells.Item(64, 1).Value2 = ${v_63_A}
$ws.Cells.Item(64, 2).Value2 = ${v_63_B}
$ws.Cells.Item(64, 4).Value2 = ${v_63_D}
$ws.Cells.Item(64, 5).Value2 = ${v_63_E}
$ws.Cells.Item(64, 6).Value2 = ${v_63_F}
$ws.Cells.Item(64, 7).Value2 = ${v_63_G}
$ws.Cells.Item(64, 8).Value2 = ${v_63_H}
$ws.Cells.Item(64, 9).Value2 = ${v_63_I}
$ws.Cells.Item(64, 16).Value2 = ${v_63_P}
$ws.Cells.Item(64, 17).Value2 = ${v_63_Q}
$ws.Cells.Item(64, 18).Value2 = ${v_63_R}
$ws.Cells.Item(64, 19).Value2 = ${v_63_S}
$ws.Cells.Item(64, 20).Value2 = ${v_63_T}
$ws.Cells.Item(64, 21).Value2 = ${v_63_U}
$ws.Cells.Item(64, 22).Value2 = ${v_63_V}
$ws.Cells.Item(64, 23).Value2 = ${v_63_W}
$ws.Cells.Item(64, 25).NumberFormat = "@"
$ws.Cells.Item(64, 25).Value2 = ${v_63_Y}
$ws.Cells.Item(64, 27).NumberFormat = "@"
$ws.Cells.Item(64, 27).Value2 = ${v_63_AA}
$ws.Cells.Item(64, 29).Value2 = ${v_63_AC}
$ws.Cells.Item(64, 30).Value2 = ${v_63_AD}
$ws.Cells.Item(64, 31).Value2 = ${v_63_AE}
$ws.Cells.Item(64, 33).Value2 = ${v_63_AG}
$ws.Cells.Item(64, 46).Value2 = ${v_63_AT}
$ws.Cells.Item(64, 49).Value2 = ${v_63_AW}
$ws.Cells.Item(64, 50).Value2 = ${v_63_AX}
$ws.Cells.Item(64, 51).Value2 = ${v_63_AY}

$ws.Cells.Item(63, 1).Value2 = ${v_64_A}
$ws.Cells.Item(63, 2).Value2 = ${v_64_B}
$ws.Cells.Item(63, 4).Value2 = ${v_64_D}
$ws.Cells.Item(63, 5).Value2 = ${v_64_E}
$ws.Cells.Item(63, 6).Value2 = ${v_64_F}
$ws.Cells.Item(63, 7).Value2 = ${v_64_G}
$ws.Cells.Item(63, 8).Value2 = ${v_64_H}
$ws.Cells.Item(63, 9).Value2 = ${v_64_I}
$ws.Cells.Item(63, 16).Value2 = ${v_64_P}
$ws.Cells.Item(63, 17).Value2 = ${v_64_Q}
$ws.Cells.Item(63, 18).Value2 = ${v_64_R}
$ws.Cells.Item(63, 19).Value2 = ${v_64_S}
$ws.Cells.Item(63, 20).Value2 = ${v_64_T}
$ws.Cells.Item(63, 21).Value2 = ${v_64_U}
$ws.Cells.Item(63, 22).Value2 = ${v_64_V}
$ws.Cells.Item(63, 23).Value2 = ${v_64_W}
$ws.Cells.Item(63, 25).NumberFormat = "@"
$ws.Cells.Item(63, 25).Value2 = ${v_64_Y}
$ws.Cells.Item(63, 27).NumberFormat = "@"
$ws.Cells.Item(63, 27).Value2 = ${v_64_AA}
$ws.Cells.Item(63, 29).Value2 = ${v_64_AC}
$ws.Cells.Item(63, 30).Value2 = ${v_64_AD}
$ws.Cells.Item(63, 31).Value2 = ${v_64_AE}
$ws.Cells.Item(63, 33).Value2 = ${v_64_AG}
$ws.Cells.Item(63, 46).Value2 = ${v_64_AT}
$ws.Cells.Item(63, 49).Value2 = ${v_64_AW}
$ws.Cells.Item(63, 50).Value2 = ${v_64_AX}
$ws.Cells.Item(63, 51).Value2 = ${v_64_AY}

# --- Rotate rows [77, 81] (content shifts forward, last wraps to first) ---
${v_77_A} = 131025492
${v_77_B} = 91804
${v_77_D} = 'NT'
${v_77_E} = 1108
${v_77_F} = 'Harticka'
${v_77_G} = 'Pelloporus leporinus'
${v_77_H} = '(Fr.) Krieglst.'
${v_77_I} = ''
${v_77_J} = ''
${v_77_K} = 'teleomorf'
${v_77_N} = ''
${v_77_P} = 'Sjöarsjön Nordöst, Jmt'
${v_77_Q} = 464400
${v_77_R} = 7042294
${v_77_S} = 10
${v_77_T} = 'Jämtland'
${v_77_U} = 'Krokom'
${v_77_V} = 'Jämtland'
${v_77_W} = 'Offerdal'
${v_77_Y} = '2026-02-03'
${v_77_AA} = '2026-02-03'
${v_77_AC} = 'Flera fruktkroppar i en stående döende gran.'
${v_77_AD} = $false
${v_77_AE} = $false
${v_77_AF} = ''
${v_77_AG} = $false
${v_77_AH} = 'Granskog'
${v_77_AJ} = 'gran'
${v_77_AK} = 'Picea abies'
${v_77_AM} = 'Trädstam på levande träd'
${v_77_AO} = 'Stem on living tree # Picea abies'
${v_77_AT} = ''
${v_77_AW} = 'Kristian Zackrisson'
${v_77_AX} = 'Kristian Zackrisson'
${v_77_AY} = ''

${v_81_A} = 131025490
${v_81_B} = 80349
${v_81_D} = 'NT'
${v_81_E} = 2081
${v_81_F} = 'Skrovellav'
${v_81_G} = 'Lobaria scrobiculata'
${v_81_H} = '(Scop.) DC.'
${v_81_I} = ''
${v_81_J} = ''
${v_81_K} = 'med soral'
${v_81_N} = ''
${v_81_P} = 'Sjöarsjön Nordöst, Jmt'
${v_81_Q} = 463994
${v_81_R} = 7042378
${v_81_S} = 10
${v_81_T} = 'Jämtland'
${v_81_U} = 'Krokom'
${v_81_V} = 'Jämtland'
${v_81_W} = 'Offerdal'
${v_81_Y} = '2026-02-03'
${v_81_AA} = '2026-02-03'
${v_81_AC} = 'Rikligt med skrovellav på en grov gammal sälg.'
${v_81_AD} = $false
${v_81_AE} = $false
${v_81_AF} = ''
${v_81_AG} = $false
${v_81_AH} = 'Granskog'
${v_81_AJ} = 'sälg'
${v_81_AK} = 'Salix caprea'
${v_81_AM} = 'Bark på levande träd'
${v_81_AO} = 'Bark on living woody plant # Salix caprea'
${v_81_AT} = ''
${v_81_AW} = 'Kristian Zackrisson'
${v_81_AX} = 'Kristian Zackrisson'
${v_81_AY} = ''

$ws.Cells.Item(81, 1).Value2 = ${v_77_A}
$ws.Cells.Item(81, 2).Value2 = ${v_77_B}
$ws.Cells.Item(81, 4).Value2 = ${v_77_D}
$ws.Cells.Item(81, 5).Value2 = ${v_77_E}
$ws.Cells.Item(81, 6).Value2 = ${v_77_F}
$ws.Cells.Item(81, 7).Value2 = ${v_77_G}
$ws.Cells.Item(81, 8).Value2 = ${v_77_H}
$ws.Cells.Item(81, 9).Value2 = ${v_77_I}
$ws.Cells.Item(81, 10).Value2 = ${v_77_J}
$ws.Cells.Item(81, 11).Value2 = ${v_77_K}
$ws.Cells.Item(81, 14).Value2 = ${v_77_N}
$ws.Cells.Item(81, 16).Value2 = ${v_77_P}
$ws.Cells.Item(81, 17).Value2 = ${v_77_Q}
$ws.Cells.Item(81, 18).Value2 = ${v_77_R}
$ws.Cells.Item(81, 19).Value2 = ${v_77_S}
$ws.Cells.Item(81, 20).Value2 = ${v_77_T}
$ws.Cells.Item(81, 21).Value2 = ${v_77_U}
$ws.Cells.Item(81, 22).Value2 = ${v_77_V}
$ws.Cells.Item(81, 23).Value2 = ${v_77_W}
$ws.Cells.Item(81, 25).NumberFormat = "@"
$ws.Cells.Item(81, 25).Value2 = ${v_77_Y}
$ws.Cells.Item(81, 27).NumberFormat = "@"
$ws.Cells.Item(81, 27).Value2 = ${v_77_AA}
$ws.Cells.Item(81, 29).Value2 = ${v_77_AC}
$ws.Cells.Item(81, 30).Value2 = ${v_77_AD}
$ws.Cells.Item(81, 31).Value2 = ${v_77_AE}
$ws.Cells.Item(81, 32).Value2 = ${v_77_AF}
$ws.Cells.Item(81, 33).Value2 = ${v_77_AG}
$ws.Cells.Item(81, 34).Value2 = ${v_77_AH}
$ws.Cells.Item(81, 36).Value2 = ${v_77_AJ}
$ws.Cells.Item(81, 37).Value2 = ${v_77_AK}
$ws.Cells.Item(81, 39).Value2 = ${v_77_AM}
$ws.Cells.Item(81, 41).Value2 = ${v_77_AO}
$ws.Cells.Item(81, 46).Value2 = ${v_77_AT}
$ws.Cells.Item(81, 49).Value2 = ${v_77_AW}
$ws.Cells.Item(81, 50).Value2 = ${v_77_AX}
$ws.Cells.Item(81, 51).Value2 = ${v_77_AY}

$ws.Cells.Item(77, 1).Value2 = ${v_81_A}
$ws.Cells.Item(77, 2).Value2 = ${v_81_B}
$ws.Cells.Item(77, 4).Value2 = ${v_81_D}
$ws.Cells.Item(77, 5).Value2 = ${v_81_E}
$ws.Cells.Item(77, 6).Value2 = ${v_81_F}
$ws.Cells.Item(77, 7).Value2 = ${v_81_G}
$ws.Cells.Item(77, 8).Value2 = ${v_81_H}
$ws.Cells.Item(77, 9).Value2 = ${v_81_I}
$ws.Cells.Item(77, 10).Value2 = ${v_81_J}
$ws.Cells.Item(77, 11).Value2 = ${v_81_K}
$ws.Cells.Item(77, 14).Value2 = ${v_81_N}
$ws.Cells.Item(77, 16).Value2 = ${v_81_P}
$ws.Cells.Item(77, 17).Value2 = ${v_81_Q}
$ws.Cells.Item(77, 18).Value2 = ${v_81_R}
$ws.Cells.Item(77, 19).Value2 = ${v_81_S}
$ws.Cells.Item(77, 20).Value2 = ${v_81_T}
$ws.Cells.Item(77, 21).Value2 = ${v_81_U}
$ws.Cells.Item(77, 22).Value2 = ${v_81_V}
$ws.Cells.Item(77, 23).Value2 = ${v_81_W}
$ws.Cells.Item(77, 25).NumberFormat = "@"
$ws.Cells.Item(77, 25).Value2 = ${v_81_Y}
$ws.Cells.Item(77, 27).NumberFormat = "@"
$ws.Cells.Item(77, 27).Value2 = ${v_81_AA}
$ws.Cells.Item(77, 29).Value2 = ${v_81_AC}
$ws.Cells.Item(77, 30).Value2 = ${v_81_AD}
$ws.Cells.Item(77, 31).Value2 = ${v_81_AE}
$ws.Cells.Item(77, 32).Value2 = ${v_81_AF}
$ws.Cells.Item(77, 33).Value2 = ${v_81_AG}
$ws.Cells.Item(77, 34).Value2 = ${v_81_AH}
$ws.Cells.Item(77, 36).Value2 = ${v_81_AJ}
$ws.Cells.Item(77, 37).Value2 = ${v_81_AK}
$ws.Cells.Item(77, 39).Value2 = ${v_81_AM}
$ws.Cells.Item(77, 41).Value2 = ${v_81_AO}
$ws.Cells.Item(77, 46).Value2 = ${v_81_AT}
$ws.Cells.Item(77, 49).Value2 = ${v_81_AW}
$ws.Cells.Item(77, 50).Value2 = ${v_81_AX}
$ws.Cells.Item(77, 51).Value2 = ${v_81_AY}

# --- Rotate rows [80, 82] (content shifts forward, last wraps to first) ---
${v_80_A} = 131024470
${v_80_B} = 57884
${v_80_D} = 'NT'
${v_80_E} = 100109
${v_80_F} = 'Tretåig hackspett'
${v_80_G} = 'Picoides tridactylus'
${v_80_H} = '(Linnaeus, 1758)'
${v_80_I} = ''
${v_80_P} = 'Sjöarsjön N, Jmt'
${v_80_Q} = 464325
${v_80_R} = 7042334
${v_80_S} = 10
${v_80_T} = 'Jämtland'
${v_80_U} = 'Krokom'
${v_80_V} = 'Jämtland'
${v_80_W} = 'Offerdal'
${v_80_Y} = '2026-02-03'
${v_80_AA} = '2026-02-03'
${v_80_AC} = 'Ringhack färska och äldre'
${v_80_AD} = $false
${v_80_AE} = $false
${v_80_AG} = $false
${v_80_AT} = ''
${v_80_AW} = 'Benny Öwre'
${v_80_AX} = 'Benny Öwre'
${v_80_AY} = ''

${v_82_A} = 131024461
${v_82_B} = 57884
${v_82_D} = 'NT'
${v_82_E} = 100109
${v_82_F} = 'Tretåig hackspett'
${v_82_G} = 'Picoides tridactylus'
${v_82_H} = '(Linnaeus, 1758)'
${v_82_I} = ''
${v_82_P} = 'Sjöarsjön N, Jmt'
${v_82_Q} = 464377
${v_82_R} = 7042402
${v_82_S} = 10
${v_82_T} = 'Jämtland'
${v_82_U} = 'Krokom'
${v_82_V} = 'Jämtland'
${v_82_W} = 'Offerdal'
${v_82_Y} = '2026-02-03'
${v_82_AA} = '2026-02-03'
${v_82_AC} = 'Ringhack äldre'
${v_82_AD} = $false
${v_82_AE} = $false
${v_82_AG} = $false
${v_82_AT} = ''
${v_82_AW} = 'Benny Öwre'
${v_82_AX} = 'Benny Öwre'
${v_82_AY} = ''

$ws.Cells.Item(82, 1).Value2 = ${v_80_A}
$ws.Cells.Item(82, 2).Value2 = ${v_80_B}
$ws.Cells.Item(82, 4).Value2 = ${v_80_D}
$ws.Cells.Item(82, 5).Value2 = ${v_80_E}
$ws.Cells.Item(82, 6).Value2 = ${v_80_F}
$ws.Cells.Item(82, 7).Value2 = ${v_80_G}
$ws.Cells.Item(82, 8).Value2 = ${v_80_H}
$ws.Cells.Item(82, 9).Value2 = ${v_80_I}
$ws.Cells.Item(82, 16).Value2 = ${v_80_P}
$ws.Cells.Item(82, 17).Value2 = ${v_80_Q}
$ws.Cells.Item(82, 18).Value2 = ${v_80_R}
$ws.Cells.Item(82, 19).Value2 = ${v_80_S}
$ws.Cells.Item(82, 20).Value2 = ${v_80_T}
$ws.Cells.Item(82, 21).Value2 = ${v_80_U}
$ws.Cells.Item(82, 22).Value2 = ${v_80_V}
$ws.Cells.Item(82, 23).Value2 = ${v_80_W}
$ws.Cells.Item(82, 25).NumberFormat = "@"
$ws.Cells.Item(82, 25).Value2 = ${v_80_Y}
$ws.Cells.Item(82, 27).NumberFormat = "@"
$ws.Cells.Item(82, 27).Value2 = ${v_80_AA}
$ws.Cells.Item(82, 29).Value2 = ${v_80_AC}
$ws.Cells.Item(82, 30).Value2 = ${v_80_AD}
$ws.Cells.Item(82, 31).Value2 = ${v_80_AE}
$ws.Cells.Item(82, 33).Value2 = ${v_80_AG}
$ws.Cells.Item(82, 46).Value2 = ${v_80_AT}
$ws.Cells.Item(82, 49).Value2 = ${v_80_AW}
$ws.Cells.Item(82, 50).Value2 = ${v_80_AX}
$ws.Cells.Item(82, 51).Value2 = ${v_80_AY}

$ws.Cells.Item(80, 1).Value2 = ${v_82_A}
$ws.Cells.Item(80, 2).Value2 = ${v_82_B}
$ws.Cells.Item(80, 4).Value2 = ${v_82_D}
$ws.Cells.Item(80, 5).Value2 = ${v_82_E}
$ws.Cells.Item(80, 6).Value2 = ${v_82_F}
$ws.Cells.Item(80, 7).Value2 = ${v_82_G}
$ws.Cells.Item(80, 8).Value2 = ${v_82_H}
$ws.Cells.Item(80, 9).Value2 = ${v_82_I}
$ws.Cells.Item(80, 16).Value2 = ${v_82_P}
$ws.Cells.Item(80, 17).Value2 = ${v_82_Q}
$ws.Cells.Item(80, 18).Value2 = ${v_82_R}
$ws.Cells.Item(80, 19).Value2 = ${v_82_S}
$ws.Cells.Item(80, 20).Value2 = ${v_82_T}
$ws.Cells.Item(80, 21).Value2 = ${v_82_U}
$ws.Cells.Item(80, 22).Value2 = ${v_82_V}
$ws.Cells.Item(80, 23).Value2 = ${v_82_W}
$ws.Cells.Item(80, 25).NumberFormat = "@"
$ws.Cells.Item(80, 25).Value2 = ${v_82_Y}
$ws.Cells.Item(80, 27).NumberFormat = "@"
$ws.Cells.Item(80, 27).Value2 = ${v_82_AA}
$ws.Cells.Item(80, 29).Value2 = ${v_82_AC}
$ws.Cells.Item(80, 30).Value2 = ${v_82_AD}
$ws.Cells.Item(80, 31).Value2 = ${v_82_AE}
$ws.Cells.Item(80, 33).Value2 = ${v_82_AG}
$ws.Cells.Item(80, 46).Value2 = ${v_82_AT}
$ws.Cells.Item(80, 49).Value2 = ${v_82_AW}
$ws.Cells.Item(80, 50).Value2 = ${v_82_AX}
$ws.Cells.Item(80, 51).Value2 = ${v_82_AY}

# --- Rotate rows [83, 84, 85, 88, 87, 86] (content shifts forward, last wraps to first) ---
${v_83_A} = 131025475
${v_83_B} = 57884
${v_83_D} = 'NT'
${v_83_E} = 100109
${v_83_F} = 'Tretåig hackspett'
${v_83_G} = 'Picoides tridactylus'
${v_83_H} = '(Linnaeus, 1758)'
${v_83_I} = ''
${v_83_K} = ''
${v_83_L} = ''
${v_83_M} = 'äldre spår'
${v_83_N} = ''
${v_83_P} = 'Sjöarsjön Nordöst, Jmt'
${v_83_Q} = 464393
${v_83_R} = 7042284
${v_83_S} = 10
${v_83_T} = 'Jämtland'
${v_83_U} = 'Krokom'
${v_83_V} = 'Jämtland'
${v_83_W} = 'Offerdal'
${v_83_Y} = '2026-02-03'
${v_83_AA} = '2026-02-03'
${v_83_AC} = 'Ringhack, äldre, på gran.'
${v_83_AD} = $false
${v_83_AE} = $false
${v_83_AG} = $false
${v_83_AH} = 'Granskog'
${v_83_AJ} = 'gran'
${v_83_AK} = 'Picea abies'
${v_83_AM} = 'Trädstam på levande träd'
${v_83_AO} = 'Stem on living tree # Picea abies'
${v_83_AT} = ''
${v_83_AW} = 'Kristian Zackrisson'
${v_83_AX} = 'Kristian Zackrisson'
${v_83_AY} = ''

${v_84_A} = 131025477
${v_84_B} = 57884
${v_84_D} = 'NT'
${v_84_E} = 100109
${v_84_F} = 'Tretåig hackspett'
${v_84_G} = 'Picoides tridactylus'
${v_84_H} = '(Linnaeus, 1758)'
${v_84_I} = ''
${v_84_K} = ''
${v_84_L} = ''
${v_84_M} = 'äldre spår'
${v_84_N} = ''
${v_84_P} = 'Sjöarsjön Nordöst, Jmt'
${v_84_Q} = 464374
${v_84_R} = 7042300
${v_84_S} = 10
${v_84_T} = 'Jämtland'
${v_84_U} = 'Krokom'
${v_84_V} = 'Jämtland'
${v_84_W} = 'Offerdal'
${v_84_Y} = '2026-02-03'
${v_84_AA} = '2026-02-03'
${v_84_AC} = 'Ringhack, äldre, på en gran.'
${v_84_AD} = $false
${v_84_AE} = $false
${v_84_AG} = $false
${v_84_AH} = 'Granskog'
${v_84_AJ} = 'gran'
${v_84_AK} = 'Picea abies'
${v_84_AM} = 'Trädstam på levande träd'
${v_84_AO} = 'Stem on living tree # Picea abies'
${v_84_AT} = ''
${v_84_AW} = 'Kristian Zackrisson'
${v_84_AX} = 'Kristian Zackrisson'
${v_84_AY} = ''

${v_85_A} = 131024473
${v_85_B} = 57884
${v_85_D} = 'NT'
${v_85_E} = 100109
${v_85_F} = 'Tretåig hackspett'
${v_85_G} = 'Picoides tridactylus'
${v_85_H} = '(Linnaeus, 1758)'
${v_85_I} = ''
${v_85_K} = $null
${v_85_L} = $null
${v_85_M} = $null
${v_85_N} = $null
${v_85_P} = 'Sjöarsjön N, Jmt'
${v_85_Q} = 464305
${v_85_R} = 7042334
${v_85_S} = 10
${v_85_T} = 'Jämtland'
${v_85_U} = 'Krokom'
${v_85_V} = 'Jämtland'
${v_85_W} = 'Offerdal'
${v_85_Y} = '2026-02-03'
${v_85_AA} = '2026-02-03'
${v_85_AC} = 'Ringhack äldre'
${v_85_AD} = $false
${v_85_AE} = $false
${v_85_AG} = $false
${v_85_AH} = $null
${v_85_AJ} = $null
${v_85_AK} = $null
${v_85_AM} = $null
${v_85_AO} = $null
${v_85_AT} = ''
${v_85_AW} = 'Benny Öwre'
${v_85_AX} = 'Benny Öwre'
${v_85_AY} = ''

${v_88_A} = 131024508
${v_88_B} = 91828
${v_88_D} = 'NT'
${v_88_E} = 5432
${v_88_F} = 'Granticka'
${v_88_G} = 'Porodaedalea chrysoloma s.lat.'
${v_88_H} = ''
${v_88_I} = ''
${v_88_K} = $null
${v_88_L} = $null
${v_88_M} = $null
${v_88_N} = $null
${v_88_P} = 'Sjöarsjön N, Jmt'
${v_88_Q} = 464083
${v_88_R} = 7042298
${v_88_S} = 10
${v_88_T} = 'Jämtland'
${v_88_U} = 'Krokom'
${v_88_V} = 'Jämtland'
${v_88_W} = 'Offerdal'
${v_88_Y} = '2026-02-03'
${v_88_AA} = '2026-02-03'
${v_88_AC} = $null
${v_88_AD} = $false
${v_88_AE} = $false
${v_88_AG} = $false
${v_88_AH} = $null
${v_88_AJ} = $null
${v_88_AK} = $null
${v_88_AM} = $null
${v_88_AO} = $null
${v_88_AT} = ''
${v_88_AW} = 'Benny Öwre'
${v_88_AX} = 'Benny Öwre'
${v_88_AY} = ''

${v_87_A} = 131024505
${v_87_B} = 91828
${v_87_D} = 'NT'
${v_87_E} = 5432
${v_87_F} = 'Granticka'
${v_87_G} = 'Porodaedalea chrysoloma s.lat.'
${v_87_H} = ''
${v_87_I} = ''
${v_87_K} = $null
${v_87_L} = $null
${v_87_M} = $null
${v_87_N} = $null
${v_87_P} = 'Sjöarsjön N, Jmt'
${v_87_Q} = 464276
${v_87_R} = 7042457
${v_87_S} = 10
${v_87_T} = 'Jämtland'
${v_87_U} = 'Krokom'
${v_87_V} = 'Jämtland'
${v_87_W} = 'Offerdal'
${v_87_Y} = '2026-02-03'
${v_87_AA} = '2026-02-03'
${v_87_AC} = $null
${v_87_AD} = $false
${v_87_AE} = $false
${v_87_AG} = $false
${v_87_AH} = $null
${v_87_AJ} = $null
${v_87_AK} = $null
${v_87_AM} = $null
${v_87_AO} = $null
${v_87_AT} = ''
${v_87_AW} = 'Benny Öwre'
${v_87_AX} = 'Benny Öwre'
${v_87_AY} = ''

${v_86_A} = 131025520
${v_86_B} = 79243
${v_86_D} = 'NT'
${v_86_E} = 6425
${v_86_F} = 'Garnlav'
${v_86_G} = 'Alectoria sarmentosa'
${v_86_H} = '(Ach.) Ach.'
${v_86_I} = ''
${v_86_K} = $null
${v_86_L} = $null
${v_86_M} = $null
${v_86_N} = $null
${v_86_P} = 'Sjöarsjön Nordöst, Jmt'
${v_86_Q} = 464027
${v_86_R} = 7042318
${v_86_S} = 10
${v_86_T} = 'Jämtland'
${v_86_U} = 'Krokom'
${v_86_V} = 'Jämtland'
${v_86_W} = 'Offerdal'
${v_86_Y} = '2026-02-03'
${v_86_AA} = '2026-02-03'
${v_86_AC} = 'På flera granar.'
${v_86_AD} = $false
${v_86_AE} = $false
${v_86_AG} = $false
${v_86_AH} = $null
${v_86_AJ} = $null
${v_86_AK} = $null
${v_86_AM} = $null
${v_86_AO} = $null
${v_86_AT} = ''
${v_86_AW} = 'Kristian Zackrisson'
${v_86_AX} = 'Kristian Zackrisson'
${v_86_AY} = ''

$ws.Cells.Item(84, 1).Value2 = ${v_83_A}
$ws.Cells.Item(84, 2).Value2 = ${v_83_B}
$ws.Cells.Item(84, 4).Value2 = ${v_83_D}
$ws.Cells.Item(84, 5).Value2 = ${v_83_E}
$ws.Cells.Item(84, 6).Value2 = ${v_83_F}
$ws.Cells.Item(84, 7).Value2 = ${v_83_G}
$ws.Cells.Item(84, 8).Value2 = ${v_83_H}
$ws.Cells.Item(84, 9).Value2 = ${v_83_I}
$ws.Cells.Item(84, 11).Value2 = ${v_83_K}
$ws.Cells.Item(84, 12).Value2 = ${v_83_L}
$ws.Cells.Item(84, 13).Value2 = ${v_83_M}
$ws.Cells.Item(84, 14).Value2 = ${v_83_N}
$ws.Cells.Item(84, 16).Value2 = ${v_83_P}
$ws.Cells.Item(84, 17).Value2 = ${v_83_Q}
$ws.Cells.Item(84, 18).Value2 = ${v_83_R}
$ws.Cells.Item(84, 19).Value2 = ${v_83_S}
$ws.Cells.Item(84, 20).Value2 = ${v_83_T}
$ws.Cells.Item(84, 21).Value2 = ${v_83_U}
$ws.Cells.Item(84, 22).Value2 = ${v_83_V}
$ws.Cells.Item(84, 23).Value2 = ${v_83_W}
$ws.Cells.Item(84, 25).NumberFormat = "@"
$ws.Cells.Item(84, 25).Value2 = ${v_83_Y}
$ws.Cells.Item(84, 27).NumberFormat = "@"
$ws.Cells.Item(84, 27).Value2 = ${v_83_AA}
$ws.Cells.Item(84, 29).Value2 = ${v_83_AC}
$ws.Cells.Item(84, 30).Value2 = ${v_83_AD}
$ws.Cells.Item(84, 31).Value2 = ${v_83_AE}
$ws.Cells.Item(84, 33).Value2 = ${v_83_AG}
$ws.Cells.Item(84, 34).Value2 = ${v_83_AH}
$ws.Cells.Item(84, 36).Value2 = ${v_83_AJ}
$ws.Cells.Item(84, 37).Value2 = ${v_83_AK}
$ws.Cells.Item(84, 39).Value2 = ${v_83_AM}
$ws.Cells.Item(84, 41).Value2 = ${v_83_AO}
$ws.Cells.Item(84, 46).Value2 = ${v_83_AT}
$ws.Cells.Item(84, 49).Value2 = ${v_83_AW}
$ws.Cells.Item(84, 50).Value2 = ${v_83_AX}
$ws.Cells.Item(84, 51).Value2 = ${v_83_AY}

$ws.Cells.Item(85, 1).Value2 = ${v_84_A}
$ws.Cells.Item(85, 2).Value2 = ${v_84_B}
$ws.Cells.Item(85, 4).Value2 = ${v_84_D}
$ws.Cells.Item(85, 5).Value2 = ${v_84_E}
$ws.Cells.Item(85, 6).Value2 = ${v_84_F}
$ws.Cells.Item(85, 7).Value2 = ${v_84_G}
$ws.Cells.Item(85, 8).Value2 = ${v_84_H}
$ws.Cells.Item(85, 9).Value2 = ${v_84_I}
$ws.Cells.Item(85, 11).Value2 = ${v_84_K}
$ws.Cells.Item(85, 12).Value2 = ${v_84_L}
$ws.Cells.Item(85, 13).Value2 = ${v_84_M}
$ws.Cells.Item(85, 14).Value2 = ${v_84_N}
$ws.Cells.Item(85, 16).Value2 = ${v_84_P}
$ws.Cells.Item(85, 17).Value2 = ${v_84_Q}
$ws.Cells.Item(85, 18).Value2 = ${v_84_R}
$ws.Cells.Item(85, 19).Value2 = ${v_84_S}
$ws.Cells.Item(85, 20).Value2 = ${v_84_T}
$ws.Cells.Item(85, 21).Value2 = ${v_84_U}
$ws.Cells.Item(85, 22).Value2 = ${v_84_V}
$ws.Cells.Item(85, 23).Value2 = ${v_84_W}
$ws.Cells.Item(85, 25).NumberFormat = "@"
$ws.Cells.Item(85, 25).Value2 = ${v_84_Y}
$ws.Cells.Item(85, 27).NumberFormat = "@"
$ws.Cells.Item(85, 27).Value2 = ${v_84_AA}
$ws.Cells.Item(85, 29).Value2 = ${v_84_AC}
$ws.Cells.Item(85, 30).Value2 = ${v_84_AD}
$ws.Cells.Item(85, 31).Value2 = ${v_84_AE}
$ws.Cells.Item(85, 33).Value2 = ${v_84_AG}
$ws.Cells.Item(85, 34).Value2 = ${v_84_AH}
$ws.Cells.Item(85, 36).Value2 = ${v_84_AJ}
$ws.Cells.Item(85, 37).Value2 = ${v_84_AK}
$ws.Cells.Item(85, 39).Value2 = ${v_84_AM}
$ws.Cells.Item(85, 41).Value2 = ${v_84_AO}
$ws.Cells.Item(85, 46).Value2 = ${v_84_AT}
$ws.Cells.Item(85, 49).Value2 = ${v_84_AW}
$ws.Cells.Item(85, 50).Value2 = ${v_84_AX}
$ws.Cells.Item(85, 51).Value2 = ${v_84_AY}

$ws.Cells.Item(88, 1).Value2 = ${v_85_A}
$ws.Cells.Item(88, 2).Value2 = ${v_85_B}
$ws.Cells.Item(88, 4).Value2 = ${v_85_D}
$ws.Cells.Item(88, 5).Value2 = ${v_85_E}
$ws.Cells.Item(88, 6).Value2 = ${v_85_F}
$ws.Cells.Item(88, 7).Value2 = ${v_85_G}
$ws.Cells.Item(88, 8).Value2 = ${v_85_H}
$ws.Cells.Item(88, 9).Value2 = ${v_85_I}
$ws.Cells.Item(88, 11).Value2 = ${v_85_K}
$ws.Cells.Item(88, 12).Value2 = ${v_85_L}
$ws.Cells.Item(88, 13).Value2 = ${v_85_M}
$ws.Cells.Item(88, 14).Value2 = ${v_85_N}
$ws.Cells.Item(88, 16).Value2 = ${v_85_P}
$ws.Cells.Item(88, 17).Value2 = ${v_85_Q}
$ws.Cells.Item(88, 18).Value2 = ${v_85_R}
$ws.Cells.Item(88, 19).Value2 = ${v_85_S}
$ws.Cells.Item(88, 20).Value2 = ${v_85_T}
$ws.Cells.Item(88, 21).Value2 = ${v_85_U}
$ws.Cells.Item(88, 22).Value2 = ${v_85_V}
$ws.Cells.Item(88, 23).Value2 = ${v_85_W}
$ws.Cells.Item(88, 25).NumberFormat = "@"
$ws.Cells.Item(88, 25).Value2 = ${v_85_Y}
$ws.Cells.Item(88, 27).NumberFormat = "@"
$ws.Cells.Item(88, 27).Value2 = ${v_85_AA}
$ws.Cells.Item(88, 29).Value2 = ${v_85_AC}
$ws.Cells.Item(88, 30).Value2 = ${v_85_AD}
$ws.Cells.Item(88, 31).Value2 = ${v_85_AE}
$ws.Cells.Item(88, 33).Value2 = ${v_85_AG}
$ws.Cells.Item(88, 34).Value2 = ${v_85_AH}
$ws.Cells.Item(88, 36).Value2 = ${v_85_AJ}
$ws.Cells.Item(88, 37).Value2 = ${v_85_AK}
$ws.Cells.Item(88, 39).Value2 = ${v_85_AM}
$ws.Cells.Item(88, 41).Value2 = ${v_85_AO}
$ws.Cells.Item(88, 46).Value2 = ${v_85_AT}
$ws.Cells.Item(88, 49).Value2 = ${v_85_AW}
$ws.Cells.Item(88, 50).Value2 = ${v_85_AX}
$ws.Cells.Item(88, 51).Value2 = ${v_85_AY}

$ws.Cells.Item(87, 1).Value2 = ${v_88_A}
$ws.Cells.Item(87, 2).Value2 = ${v_88_B}
$ws.Cells.Item(87, 4).Value2 = ${v_88_D}
$ws.Cells.Item(87, 5).Value2 = ${v_88_E}
$ws.Cells.Item(87, 6).Value2 = ${v_88_F}
$ws.Cells.Item(87, 7).Value2 = ${v_88_G}
$ws.Cells.Item(87, 8).Value2 = ${v_88_H}
$ws.Cells.Item(87, 9).Value2 = ${v_88_I}
$ws.Cells.Item(87, 11).Value2 = ${v_88_K}
$ws.Cells.Item(87, 12).Value2 = ${v_88_L}
$ws.Cells.Item(87, 13).Value2 = ${v_88_M}
$ws.Cells.Item(87, 14).Value2 = ${v_88_N}
$ws.Cells.Item(87, 16).Value2 = ${v_88_P}
$ws.Cells.Item(87, 17).Value2 = ${v_88_Q}
$ws.Cells.Item(87, 18).Value2 = ${v_88_R}
$ws.Cells.Item(87, 19).Value2 = ${v_88_S}
$ws.Cells.Item(87, 20).Value2 = ${v_88_T}
$ws.Cells.Item(87, 21).Value2 = ${v_88_U}
$ws.Cells.Item(87, 22).Value2 = ${v_88_V}
$ws.Cells.Item(87, 23).Value2 = ${v_88_W}
$ws.Cells.Item(87, 25).NumberFormat = "@"
$ws.Cells.Item(87, 25).Value2 = ${v_88_Y}
$ws.Cells.Item(87, 27).NumberFormat = "@"
$ws.Cells.Item(87, 27).Value2 = ${v_88_AA}
$ws.Cells.Item(87, 29).Value2 = ${v_88_AC}
$ws.Cells.Item(87, 30).Value2 = ${v_88_AD}
$ws.Cells.Item(87, 31).Value2 = ${v_88_AE}
$ws.Cells.Item(87, 33).Value2 = ${v_88_AG}
$ws.Cells.Item(87, 34).Value2 = ${v_88_AH}
$ws.Cells.Item(87, 36).Value2 = ${v_88_AJ}
$ws.Cells.Item(87, 37).Value2 = ${v_88_AK}
$ws.Cells.Item(87, 39).Value2 = ${v_88_AM}
$ws.Cells.Item(87, 41).Value2 = ${v_88_AO}
$ws.Cells.Item(87, 46).Value2 = ${v_88_AT}
$ws.Cells.Item(87, 49).Value2 = ${v_88_AW}
$ws.Cells.Item(87, 50).Value2 = ${v_88_AX}
$ws.Cells.Item(87, 51).Value2 = ${v_88_AY}

$ws.Cells.Item(86, 1).Value2 = ${v_87_A}
$ws.Cells.Item(86, 2).Value2 = ${v_87_B}
$ws.Cells.Item(86, 4).Value2 = ${v_87_D}
$ws.Cells.Item(86, 5).Value2 = ${v_87_E}
$ws.Cells.Item(86, 6).Value2 = ${v_87_F}
$ws.Cells.Item(86, 7).Value2 = ${v_87_G}
$ws.Cells.Item(86, 8).Value2 = ${v_87_H}
$ws.Cells.Item(86, 9).Value2 = ${v_87_I}
$ws.Cells.Item(86, 11).Value2 = ${v_87_K}
$ws.Cells.Item(86, 12).Value2 = ${v_87_L}
$ws.Cells.Item(86, 13).Value2 = ${v_87_M}
$ws.Cells.Item(86, 14).Value2 = ${v_87_N}
$ws.Cells.Item(86, 16).Value2 = ${v_87_P}
$ws.Cells.Item(86, 17).Value2 = ${v_87_Q}
$ws.Cells.Item(86, 18).Value2 = ${v_87_R}
$ws.Cells.Item(86, 19).Value2 = ${v_87_S}
$ws.Cells.Item(86, 20).Value2 = ${v_87_T}
$ws.Cells.Item(86, 21).Value2 = ${v_87_U}
$ws.Cells.Item(86, 22).Value2 = ${v_87_V}
$ws.Cells.Item(86, 23).Value2 = ${v_87_W}
$ws.Cells.Item(86, 25).NumberFormat = "@"
$ws.Cells.Item(86, 25).Value2 = ${v_87_Y}
$ws.Cells.Item(86, 27).NumberFormat = "@"
$ws.Cells.Item(86, 27).Value2 = ${v_87_AA}
$ws.Cells.Item(86, 29).Value2 = ${v_87_AC}
$ws.Cells.Item(86, 30).Value2 = ${v_87_AD}
$ws.Cells.Item(86, 31).Value2 = ${v_87_AE}
$ws.Cells.Item(86, 33).Value2 = ${v_87_AG}
$ws.Cells.Item(86, 34).Value2 = ${v_87_AH}
$ws.Cells.Item(86, 36).Value2 = ${v_87_AJ}
$ws.Cells.Item(86, 37).Value2 = ${v_87_AK}
$ws.Cells.Item(86, 39).Value2 = ${v_87_AM}
$ws.Cells.Item(86, 41).Value2 = ${v_87_AO}
$ws.Cells.Item(86, 46).Value2 = ${v_87_AT}
$ws.Cells.Item(86, 49).Value2 = ${v_87_AW}
$ws.Cells.Item(86, 50).Value2 = ${v_87_AX}
$ws.Cells.Item(86, 51).Value2 = ${v_87_AY}

$ws.Cells.Item(83, 1).Value2 = ${v_86_A}
$ws.Cells.Item(83, 2).Value2 = ${v_86_B}
$ws.Cells.Item(83, 4).Value2 = ${v_86_D}
$ws.Cells.Item(83, 5).Value2 = ${v_86_E}
$ws.Cells.Item(83, 6).Value2 = ${v_86_F}
$ws.Cells.Item(83, 7).Value2 = ${v_86_G}
$ws.Cells.Item(83, 8).Value2 = ${v_86_H}
$ws.Cells.Item(83, 9).Value2 = ${v_86_I}
$ws.Cells.Item(83, 11).Value2 = ${v_86_K}
$ws.Cells.Item(83, 12).Value2 = ${v_86_L}
$ws.Cells.Item(83, 13).Value2 = ${v_86_M}
$ws.Cells.Item(83, 14).Value2 = ${v_86_N}
$ws.Cells.Item(83, 16).Value2 = ${v_86_P}
$ws.Cells.Item(83, 17).Value2 = ${v_86_Q}
$ws.Cells.Item(83, 18).Value2 = ${v_86_R}
$ws.Cells.Item(83, 19).Value2 = ${v_86_S}
$ws.Cells.Item(83, 20).Value2 = ${v_86_T}
$ws.Cells.Item(83, 21).Value2 = ${v_86_U}
$ws.Cells.Item(83, 22).Value2 = ${v_86_V}
$ws.Cells.Item(83, 23).Value2 = ${v_86_W}
$ws.Cells.Item(83, 25).NumberFormat = "@"
$ws.Cells.Item(83, 25).Value2 = ${v_86_Y}
$ws.Cells.Item(83, 27).NumberFormat = "@"
$ws.Cells.Item(83, 27).Value2 = ${v_86_AA}
$ws.Cells.Item(83, 29).Value2 = ${v_86_AC}
$ws.Cells.Item(83, 30).Value2 = ${v_86_AD}
$ws.Cells.Item(83, 31).Value2 = ${v_86_AE}
$ws.Cells.Item(83, 33).Value2 = ${v_86_AG}
$ws.Cells.Item(83, 34).Value2 = ${v_86_AH}
$ws.Cells.Item(83, 36).Value2 = ${v_86_AJ}
$ws.Cells.Item(83, 37).Value2 = ${v_86_AK}
$ws.Cells.Item(83, 39).Value2 = ${v_86_AM}
$ws.Cells.Item(83, 41).Value2 = ${v_86_AO}
$ws.Cells.Item(83, 46).Value2 = ${v_86_AT}
$ws.Cells.Item(83, 49).Value2 = ${v_86_AW}
$ws.Cells.Item(83, 50).Value2 = ${v_86_AX}
$ws.Cells.Item(83, 51).Value2 = ${v_86_AY}

# --- Rotate rows [91, 92, 93] (content shifts forward, last wraps to first) ---
${v_91_A} = 131024455
${v_91_B} = 57884
${v_91_D} = 'NT'
${v_91_E} = 100109
${v_91_F} = 'Tretåig hackspett'
${v_91_G} = 'Picoides tridactylus'
${v_91_H} = '(Linnaeus, 1758)'
${v_91_I} = ''
${v_91_P} = 'Sjöarsjön N, Jmt'
${v_91_Q} = 464454
${v_91_R} = 7042523
${v_91_S} = 10
${v_91_T} = 'Jämtland'
${v_91_U} = 'Krokom'
${v_91_V} = 'Jämtland'
${v_91_W} = 'Offerdal'
${v_91_Y} = '2026-02-03'
${v_91_AA} = '2026-02-03'
${v_91_AC} = 'Ringhack äldre'
${v_91_AD} = $false
${v_91_AE} = $false
${v_91_AG} = $false
${v_91_AT} = ''
${v_91_AW} = 'Benny Öwre'
${v_91_AX} = 'Benny Öwre'
${v_91_AY} = ''

${v_92_A} = 131024466
${v_92_B} = 57884
${v_92_D} = 'NT'
${v_92_E} = 100109
${v_92_F} = 'Tretåig hackspett'
${v_92_G} = 'Picoides tridactylus'
${v_92_H} = '(Linnaeus, 1758)'
${v_92_I} = ''
${v_92_P} = 'Sjöarsjön N, Jmt'
${v_92_Q} = 464319
${v_92_R} = 7042385
${v_92_S} = 10
${v_92_T} = 'Jämtland'
${v_92_U} = 'Krokom'
${v_92_V} = 'Jämtland'
${v_92_W} = 'Offerdal'
${v_92_Y} = '2026-02-03'
${v_92_AA} = '2026-02-03'
${v_92_AC} = 'Ringhack'
${v_92_AD} = $false
${v_92_AE} = $false
${v_92_AG} = $false
${v_92_AT} = ''
${v_92_AW} = 'Benny Öwre'
${v_92_AX} = 'Benny Öwre'
${v_92_AY} = ''

${v_93_A} = 131024475
${v_93_B} = 57884
${v_93_D} = 'NT'
${v_93_E} = 100109
${v_93_F} = 'Tretåig hackspett'
${v_93_G} = 'Picoides tridactylus'
${v_93_H} = '(Linnaeus, 1758)'
${v_93_I} = ''
${v_93_P} = 'Sjöarsjön N, Jmt'
${v_93_Q} = 464272
${v_93_R} = 7042273
${v_93_S} = 10
${v_93_T} = 'Jämtland'
${v_93_U} = 'Krokom'
${v_93_V} = 'Jämtland'
${v_93_W} = 'Offerdal'
${v_93_Y} = '2026-02-03'
${v_93_AA} = '2026-02-03'
${v_93_AC} = 'Ringhack äldre'
${v_93_AD} = $false
${v_93_AE} = $false
${v_93_AG} = $false
${v_93_AT} = ''
${v_93_AW} = 'Benny Öwre'
${v_93_AX} = 'Benny Öwre'
${v_93_AY} = ''

$ws.Cells.Item(92, 1).Value2 = ${v_91_A}
$ws.Cells.Item(92, 2).Value2 = ${v_91_B}
$ws.Cells.Item(92, 4).Value2 = ${v_91_D}
$ws.Cells.Item(92, 5).Value2 = ${v_91_E}
$ws.Cells.Item(92, 6).Value2 = ${v_91_F}
$ws.Cells.Item(92, 7).Value2 = ${v_91_G}
$ws.Cells.Item(92, 8).Value2 = ${v_91_H}
$ws.Cells.Item(92, 9).Value2 = ${v_91_I}
$ws.Cells.Item(92, 16).Value2 = ${v_91_P}
$ws.Cells.Item(92, 17).Value2 = ${v_91_Q}
$ws.Cells.Item(92, 18).Value2 = ${v_91_R}
$ws.Cells.Item(92, 19).Value2 = ${v_91_S}
$ws.Cells.Item(92, 20).Value2 = ${v_91_T}
$ws.Cells.Item(92, 21).Value2 = ${v_91_U}
$ws.Cells.Item(92, 22).Value2 = ${v_91_V}
$ws.Cells.Item(92, 23).Value2 = ${v_91_W}
$ws.Cells.Item(92, 25).NumberFormat = "@"
$ws.Cells.Item(92, 25).Value2 = ${v_91_Y}
$ws.Cells.Item(92, 27).NumberFormat = "@"
$ws.Cells.Item(92, 27).Value2 = ${v_91_AA}
$ws.Cells.Item(92, 29).Value2 = ${v_91_AC}
$ws.Cells.Item(92, 30).Value2 = ${v_91_AD}
$ws.Cells.Item(92, 31).Value2 = ${v_91_AE}
$ws.Cells.Item(92, 33).Value2 = ${v_91_AG}
$ws.Cells.Item(92, 46).Value2 = ${v_91_AT}
$ws.Cells.Item(92, 49).Value2 = ${v_91_AW}
$ws.Cells.Item(92, 50).Value2 = ${v_91_AX}
$ws.Cells.Item(92, 51).Value2 = ${v_91_AY}

$ws.Cells.Item(93, 1).Value2 = ${v_92_A}
$ws.Cells.Item(93, 2).Value2 = ${v_92_B}
$ws.Cells.Item(93, 4).Value2 = ${v_92_D}
$ws.Cells.Item(93, 5).Value2 = ${v_92_E}
$ws.Cells.Item(93, 6).Value2 = ${v_92_F}
$ws.Cells.Item(93, 7).Value2 = ${v_92_G}
$ws.Cells.Item(93, 8).Value2 = ${v_92_H}
$ws.Cells.Item(93, 9).Value2 = ${v_92_I}
$ws.Cells.Item(93, 16).Value2 = ${v_92_P}
$ws.Cells.Item(93, 17).Value2 = ${v_92_Q}
$ws.Cells.Item(93, 18).Value2 = ${v_92_R}
$ws.Cells.Item(93, 19).Value2 = ${v_92_S}
$ws.Cells.Item(93, 20).Value2 = ${v_92_T}
$ws.Cells.Item(93, 21).Value2 = ${v_92_U}
$ws.Cells.Item(93, 22).Value2 = ${v_92_V}
$ws.Cells.Item(93, 23).Value2 = ${v_92_W}
$ws.Cells.Item(93, 25).NumberFormat = "@"
$ws.Cells.Item(93, 25).Value2 = ${v_92_Y}
$ws.Cells.Item(93, 27).NumberFormat = "@"
$ws.Cells.Item(93, 27).Value2 = ${v_92_AA}
$ws.Cells.Item(93, 29).Value2 = ${v_92_AC}
$ws.Cells.Item(93, 30).Value2 = ${v_92_AD}
$ws.Cells.Item(93, 31).Value2 = ${v_92_AE}
$ws.Cells.Item(93, 33).Value2 = ${v_92_AG}
$ws.Cells.Item(93, 46).Value2 = ${v_92_AT}
$ws.Cells.Item(93, 49).Value2 = ${v_92_AW}
$ws.Cells.Item(93, 50).Value2 = ${v_92_AX}
$ws.Cells.Item(93, 51).Value2 = ${v_92_AY}

$ws.Cells.Item(91, 1).Value2 = ${v_93_A}
$ws.Cells.Item(91, 2).Value2 = ${v_93_B}
$ws.Cells.Item(91, 4).Value2 = ${v_93_D}
$ws.Cells.Item(91, 5).Value2 = ${v_93_E}
$ws.Cells.Item(91, 6).Value2 = ${v_93_F}
$ws.Cells.Item(91, 7).Value2 = ${v_93_G}
$ws.Cells.Item(91, 8).Value2 = ${v_93_H}
$ws.Cells.Item(91, 9).Value2 = ${v_93_I}
$ws.Cells.Item(91, 16).Value2 = ${v_93_P}
$ws.Cells.Item(91, 17).Value2 = ${v_93_Q}
$ws.Cells.Item(91, 18).Value2 = ${v_93_R}
$ws.Cells.Item(91, 19).Value2 = ${v_93_S}
$ws.Cells.Item(91, 20).Value2 = ${v_93_T}
$ws.Cells.Item(91, 21).Value2 = ${v_93_U}
$ws.Cells.Item(91, 22).Value2 = ${v_93_V}
$ws.Cells.Item(91, 23).Value2 = ${v_93_W}
$ws.Cells.Item(91, 25).NumberFormat = "@"
$ws.Cells.Item(91, 25).Value2 = ${v_93_Y}
$ws.Cells.Item(91, 27).NumberFormat = "@"
$ws.Cells.Item(91, 27).Value2 = ${v_93_AA}
$ws.Cells.Item(91, 29).Value2 = ${v_93_AC}
$ws.Cells.Item(91, 30).Value2 = ${v_93_AD}
$ws.Cells.Item(91, 31).Value2 = ${v_93_AE}
$ws.Cells.Item(91, 33).Value2 = ${v_93_AG}
$ws.Cells.Item(91, 46).Value2 = ${v_93_AT}
$ws.Cells.Item(91, 49).Value2 = ${v_93_AW}
$ws.Cells.Item(91, 50).Value2 = ${v_93_AX}
$ws.Cells.Item(91, 51).Value2 = ${v_93_AY}

# --- Rotate rows [94, 95] (content shifts forward, last wraps to first) ---
${v_94_A} = 131025473
${v_94_B} = 57884
${v_94_D} = 'NT'
${v_94_E} = 100109
${v_94_F} = 'Tretåig hackspett'
${v_94_G} = 'Picoides tridactylus'
${v_94_H} = '(Linnaeus, 1758)'
${v_94_I} = ''
${v_94_K} = ''
${v_94_L} = ''
${v_94_M} = 'färska spår'
${v_94_N} = ''
${v_94_P} = 'Sjöarsjön Nordöst, Jmt'
${v_94_Q} = 464413
${v_94_R} = 7042312
${v_94_S} = 10
${v_94_T} = 'Jämtland'
${v_94_U} = 'Krokom'
${v_94_V} = 'Jämtland'
${v_94_W} = 'Offerdal'
${v_94_Y} = '2026-02-03'
${v_94_AA} = '2026-02-03'
${v_94_AC} = 'Ringhack, färska och äldre, längs ett par meter på en gran med spår av rikligt kåda-flöde.'
${v_94_AD} = $false
${v_94_AE} = $false
${v_94_AG} = $false
${v_94_AH} = 'Granskog'
${v_94_AI} = 'Gammal granskog.'
${v_94_AJ} = 'gran'
${v_94_AK} = 'Picea abies'
${v_94_AM} = 'Trädstam på levande träd'
${v_94_AO} = 'Stem on living tree # Picea abies'
${v_94_AT} = ''
${v_94_AW} = 'Kristian Zackrisson'
${v_94_AX} = 'Kristian Zackrisson'
${v_94_AY} = ''

${v_95_A} = 131024512
${v_95_B} = 91828
${v_95_D} = 'NT'
${v_95_E} = 5432
${v_95_F} = 'Granticka'
${v_95_G} = 'Porodaedalea chrysoloma s.lat.'
${v_95_H} = ''
${v_95_I} = ''
${v_95_K} = $null
${v_95_L} = $null
${v_95_M} = $null
${v_95_N} = $null
${v_95_P} = 'Sjöarsjön N, Jmt'
${v_95_Q} = 464059
${v_95_R} = 7042514
${v_95_S} = 10
${v_95_T} = 'Jämtland'
${v_95_U} = 'Krokom'
${v_95_V} = 'Jämtland'
${v_95_W} = 'Offerdal'
${v_95_Y} = '2026-02-03'
${v_95_AA} = '2026-02-03'
${v_95_AC} = $null
${v_95_AD} = $false
${v_95_AE} = $false
${v_95_AG} = $false
${v_95_AH} = $null
${v_95_AI} = $null
${v_95_AJ} = $null
${v_95_AK} = $null
${v_95_AM} = $null
${v_95_AO} = $null
${v_95_AT} = ''
${v_95_AW} = 'Benny Öwre'
${v_95_AX} = 'Benny Öwre'
${v_95_AY} = ''

$ws.Cells.Item(95, 1).Value2 = ${v_94_A}
$ws.Cells.Item(95, 2).Value2 = ${v_94_B}
$ws.Cells.Item(95, 4).Value2 = ${v_94_D}
$ws.Cells.Item(95, 5).Value2 = ${v_94_E}
$ws.Cells.Item(95, 6).Value2 = ${v_94_F}
$ws.Cells.Item(95, 7).Value2 = ${v_94_G}
$ws.Cells.Item(95, 8).Value2 = ${v_94_H}
$ws.Cells.Item(95, 9).Value2 = ${v_94_I}
$ws.Cells.Item(95, 11).Value2 = ${v_94_K}
$ws.Cells.Item(95, 12).Value2 = ${v_94_L}
$ws.Cells.Item(95, 13).Value2 = ${v_94_M}
$ws.Cells.Item(95, 14).Value2 = ${v_94_N}
$ws.Cells.Item(95, 16).Value2 = ${v_94_P}
$ws.Cells.Item(95, 17).Value2 = ${v_94_Q}
$ws.Cells.Item(95, 18).Value2 = ${v_94_R}
$ws.Cells.Item(95, 19).Value2 = ${v_94_S}
$ws.Cells.Item(95, 20).Value2 = ${v_94_T}
$ws.Cells.Item(95, 21).Value2 = ${v_94_U}
$ws.Cells.Item(95, 22).Value2 = ${v_94_V}
$ws.Cells.Item(95, 23).Value2 = ${v_94_W}
$ws.Cells.Item(95, 25).NumberFormat = "@"
$ws.Cells.Item(95, 25).Value2 = ${v_94_Y}
$ws.Cells.Item(95, 27).NumberFormat = "@"
$ws.Cells.Item(95, 27).Value2 = ${v_94_AA}
$ws.Cells.Item(95, 29).Value2 = ${v_94_AC}
$ws.Cells.Item(95, 30).Value2 = ${v_94_AD}
$ws.Cells.Item(95, 31).Value2 = ${v_94_AE}
$ws.Cells.Item(95, 33).Value2 = ${v_94_AG}
$ws.Cells.Item(95, 34).Value2 = ${v_94_AH}
$ws.Cells.Item(95, 35).Value2 = ${v_94_AI}
$ws.Cells.Item(95, 36).Value2 = ${v_94_AJ}
$ws.Cells.Item(95, 37).Value2 = ${v_94_AK}
$ws.Cells.Item(95, 39).Value2 = ${v_94_AM}
$ws.Cells.Item(95, 41).Value2 = ${v_94_AO}
$ws.Cells.Item(95, 46).Value2 = ${v_94_AT}
$ws.Cells.Item(95, 49).Value2 = ${v_94_AW}
$ws.Cells.Item(95, 50).Value2 = ${v_94_AX}
$ws.Cells.Item(95, 51).Value2 = ${v_94_AY}

$ws.Cells.Item(94, 1).Value2 = ${v_95_A}
$ws.Cells.Item(94, 2).Value2 = ${v_95_B}
$ws.Cells.Item(94, 4).Value2 = ${v_95_D}
$ws.Cells.Item(94, 5).Value2 = ${v_95_E}
$ws.Cells.Item(94, 6).Value2 = ${v_95_F}
$ws.Cells.Item(94, 7).Value2 = ${v_95_G}
$ws.Cells.Item(94, 8).Value2 = ${v_95_H}
$ws.Cells.Item(94, 9).Value2 = ${v_95_I}
$ws.Cells.Item(94, 11).Value2 = ${v_95_K}
$ws.Cells.Item(94, 12).Value2 = ${v_95_L}
$ws.Cells.Item(94, 13).Value2 = ${v_95_M}
$ws.Cells.Item(94, 14).Value2 = ${v_95_N}
$ws.Cells.Item(94, 16).Value2 = ${v_95_P}
$ws.Cells.Item(94, 17).Value2 = ${v_95_Q}
$ws.Cells.Item(94, 18).Value2 = ${v_95_R}
$ws.Cells.Item(94, 19).Value2 = ${v_95_S}
$ws.Cells.Item(94, 20).Value2 = ${v_95_T}
$ws.Cells.Item(94, 21).Value2 = ${v_95_U}
$ws.Cells.Item(94, 22).Value2 = ${v_95_V}
$ws.Cells.Item(94, 23).Value2 = ${v_95_W}
$ws.Cells.Item(94, 25).NumberFormat = "@"
$ws.Cells.Item(94, 25).Value2 = ${v_95_Y}
$ws.Cells.Item(94, 27).NumberFormat = "@"
$ws.Cells.Item(94, 27).Value2 = ${v_95_AA}
$ws.Cells.Item(94, 29).Value2 = ${v_95_AC}
$ws.Cells.Item(94, 30).Value2 = ${v_95_AD}
$ws.Cells.Item(94, 31).Value2 = ${v_95_AE}
$ws.Cells.Item(94, 33).Value2 = ${v_95_AG}
$ws.Cells.Item(94, 34).Value2 = ${v_95_AH}
$ws.Cells.Item(94, 35).Value2 = ${v_95_AI}
$ws.Cells.Item(94, 36).Value2 = ${v_95_AJ}
$ws.Cells.Item(94, 37).Value2 = ${v_95_AK}
$ws.Cells.Item(94, 39).Value2 = ${v_95_AM}
$ws.Cells.Item(94, 41).Value2 = ${v_95_AO}
$ws.Cells.Item(94, 46).Value2 = ${v_95_AT}
$ws.Cells.Item(94, 49).Value2 = ${v_95_AW}
$ws.Cells.Item(94, 50).Value2 = ${v_95_AX}
$ws.Cells.Item(94, 51).Value2 = ${v_95_AY}

